# Updated cryptos list (mirrors the periodic GitHub Actions refresh of
# coinranking.com data). Only the Price (D) and Volume(1h) (E) columns
# change for most rows; rows 41/42 (Filecoin <-> dogwifhat) swap order.
#
# Price values that look like plain decimals (e.g. "590.24") would be
# auto-converted to numbers by Excel if assigned directly, losing the
# original text formatting. Prefixing with a leading apostrophe forces
# Excel to keep/treat them as text, same as typing '590.24 into a cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2";  Value = "68.812.04" }
    @{ Cell = "E2";  Value = "  +2.10%  " }
    @{ Cell = "D3";  Value = "3.308.07" }
    @{ Cell = "E3";  Value = "  +2.08%  " }
    @{ Cell = "E4";  Value = "  +0.03%  " }
    @{ Cell = "D5";  Value = "590.24"; ForceText = $true }
    @{ Cell = "E5";  Value = "  +2.26%  " }
    @{ Cell = "D6";  Value = "186.89"; ForceText = $true }
    @{ Cell = "E6";  Value = "  +4.47%  " }
    @{ Cell = "E7";  Value = "  +0.01%  " }
    @{ Cell = "E8";  Value = "  +0.78%  " }
    @{ Cell = "E9";  Value = "  +4.94%  " }
    @{ Cell = "E10"; Value = "  -0.39%  " }
    @{ Cell = "E11"; Value = "  +2.72%  " }
    @{ Cell = "D12"; Value = "3.881.23" }
    @{ Cell = "E12"; Value = "  +2.19%  " }
    @{ Cell = "E13"; Value = "  +0.32%  " }
    @{ Cell = "D14"; Value = "29.06"; ForceText = $true }
    @{ Cell = "E14"; Value = "  +4.41%  " }
    @{ Cell = "D15"; Value = "68.808.24" }
    @{ Cell = "E15"; Value = "  +2.31%  " }
    @{ Cell = "E16"; Value = "  +3.74%  " }
    @{ Cell = "D17"; Value = "3.308.12" }
    @{ Cell = "E17"; Value = "  +2.23%  " }
    @{ Cell = "D18"; Value = "5.91"; ForceText = $true }
    @{ Cell = "E18"; Value = "  +2.00%  " }
    @{ Cell = "E19"; Value = "  +3.24%  " }
    @{ Cell = "D20"; Value = "385.87"; ForceText = $true }
    @{ Cell = "E20"; Value = "  +3.11%  " }
    @{ Cell = "D21"; Value = "7.83"; ForceText = $true }
    @{ Cell = "E21"; Value = "  +3.39%  " }
    @{ Cell = "D22"; Value = "71.74"; ForceText = $true }
    @{ Cell = "E22"; Value = "  +0.95%  " }
    @{ Cell = "D23"; Value = "0.999"; ForceText = $true }
    @{ Cell = "E23"; Value = "  -0.28%  " }
    @{ Cell = "E24"; Value = "  +4.00%  " }
    @{ Cell = "E25"; Value = "  +2.30%  " }
    @{ Cell = "E26"; Value = "  +7.74%  " }
    @{ Cell = "E27"; Value = "  +2.97%  " }
    @{ Cell = "D28"; Value = "1.00"; ForceText = $true }
    @{ Cell = "E28"; Value = "  -0.39%  " }
    @{ Cell = "E29"; Value = "  +5.68%  " }
    @{ Cell = "E30"; Value = "  +2.42%  " }
    @{ Cell = "E31"; Value = "  +5.60%  " }
    @{ Cell = "D32"; Value = "23.14"; ForceText = $true }
    @{ Cell = "E32"; Value = "  +2.51%  " }
    @{ Cell = "E33"; Value = "  +6.81%  " }
    @{ Cell = "E34"; Value = "  +0.01%  " }
    @{ Cell = "D35"; Value = "1.55"; ForceText = $true }
    @{ Cell = "E35"; Value = "  +4.50%  " }
    @{ Cell = "D36"; Value = "163.54"; ForceText = $true }
    @{ Cell = "E36"; Value = "  -0.19%  " }
    @{ Cell = "E37"; Value = "  +2.74%  " }
    @{ Cell = "D38"; Value = "0.841"; ForceText = $true }
    @{ Cell = "E38"; Value = "  -2.41%  " }
    @{ Cell = "D39"; Value = "27.02"; ForceText = $true }
    @{ Cell = "E39"; Value = "  +1.03%  " }
    @{ Cell = "D40"; Value = "6.78"; ForceText = $true }
    @{ Cell = "E40"; Value = "  -0.81%  " }
    @{ Cell = "B41"; Value = "dogwifhat" }
    @{ Cell = "C41"; Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif" }
    @{ Cell = "D41"; Value = "2.67"; ForceText = $true }
    @{ Cell = "E41"; Value = "  +3.88%  " }
    @{ Cell = "B42"; Value = "Filecoin" }
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil" }
    @{ Cell = "D42"; Value = "4.64"; ForceText = $true }
    @{ Cell = "E42"; Value = "  +5.72%  " }
    @{ Cell = "D43"; Value = "25.88"; ForceText = $true }
    @{ Cell = "E43"; Value = "  +0.79%  " }
    @{ Cell = "E44"; Value = "  +3.93%  " }
    @{ Cell = "D45"; Value = "41.41"; ForceText = $true }
    @{ Cell = "E45"; Value = "  +2.55%  " }
    @{ Cell = "D46"; Value = "2.653.00" }
    @{ Cell = "E46"; Value = "  -2.04%  " }
    @{ Cell = "D47"; Value = "342.67"; ForceText = $true }
    @{ Cell = "E47"; Value = "  -5.55%  " }
    @{ Cell = "E48"; Value = "  +3.37%  " }
    @{ Cell = "D49"; Value = "32.45"; ForceText = $true }
    @{ Cell = "E49"; Value = "  +6.09%  " }
    @{ Cell = "E50"; Value = "  +1.87%  " }
    @{ Cell = "E51"; Value = "  +0.52%  " }
)

foreach ($u in $updates) {
    if ($u.ForceText) {
        $ws.Range($u.Cell).Value = "'" + $u.Value
    } else {
        $ws.Range($u.Cell).Value = $u.Value
    }
}
